# Scheduled-runner update: refresh cached market-board pricing figures
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets of the Bahamut_Profits
# workbook. Only raw numeric <v> cells are touched; no formulas exist in
# this workbook, so plain Range.Value assignment is sufficient. A couple
# of rows had their (now-empty) LeveProfitHQ / LeveProfitNQ cell removed
# entirely as part of the refresh, so ClearContents() is used there.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 250183
$ws.Range("I33").Value = 333513.34
$ws.Range("J33").Value = 192
$ws.Range("K33").Value = 333513.34
$ws.Range("L33").Value = 192
$ws.Range("M33").Value = -333284.34
$ws.Range("N33").Value = -650
$ws.Range("H40").Value = 33335636
$ws.Range("J40").Value = 37039372
$ws.Range("L40").Value = 37039372
$ws.Range("N40").Value = -37039722
$ws.Range("H64").Value = 15000
$ws.Range("I64").Value = 15000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 15000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -14752
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 15000
$ws.Range("I67").Value = 15000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -14142
$ws.Range("N67").ClearContents()
$ws.Range("H118").Value = 391.42856
$ws.Range("I118").Value = 290
$ws.Range("K118").Value = 870
$ws.Range("M118").Value = 787
$ws.Range("H129").Value = 1041.3871
$ws.Range("I129").Value = 296.33334
$ws.Range("J129").Value = 1220.2
$ws.Range("K129").Value = 889.0000200000001
$ws.Range("L129").Value = 3660.6
$ws.Range("M129").Value = 4110.99998
$ws.Range("N129").Value = -13660.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2628.8684
$ws.Range("I2").Value = 2560.1428
$ws.Range("J2").Value = 2821.3
$ws.Range("K2").Value = 2560.1428
$ws.Range("L2").Value = 2821.3
$ws.Range("M2").Value = -2447.1428
$ws.Range("N2").Value = -3047.3
$ws.Range("H61").Value = 1439.3055
$ws.Range("I61").Value = 1266.4814
$ws.Range("J61").Value = 1957.7778
$ws.Range("K61").Value = 1266.4814
$ws.Range("L61").Value = 1957.7778
$ws.Range("M61").Value = -1054.4814
$ws.Range("N61").Value = -2381.7778
$ws.Range("H88").Value = 2600
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 2700
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 2700
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -3512
$ws.Range("H91").Value = 2600
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 2700
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 2700
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -5508
$ws.Range("H97").Value = 726.06665
$ws.Range("I97").Value = 711.1
$ws.Range("J97").Value = 756
$ws.Range("K97").Value = 711.1
$ws.Range("L97").Value = 756
$ws.Range("M97").Value = -215.1
$ws.Range("N97").Value = -1748
$ws.Range("H116").Value = 2628.8684
$ws.Range("I116").Value = 2560.1428
$ws.Range("J116").Value = 2821.3
$ws.Range("K116").Value = 2560.1428
$ws.Range("L116").Value = 2821.3
$ws.Range("M116").Value = -266.1428000000001
$ws.Range("N116").Value = -7409.3
$ws.Range("H132").Value = 1791.7307
$ws.Range("I132").Value = 1286.3334
$ws.Range("K132").Value = 3859.0002
$ws.Range("M132").Value = -1329.0002
$ws.Range("H136").Value = 1439.3055
$ws.Range("I136").Value = 1266.4814
$ws.Range("J136").Value = 1957.7778
$ws.Range("K136").Value = 3799.4442
$ws.Range("L136").Value = 5873.3334
$ws.Range("M136").Value = -1249.4442
$ws.Range("N136").Value = -10973.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2628.8684
$ws.Range("I3").Value = 2560.1428
$ws.Range("J3").Value = 2821.3
$ws.Range("K3").Value = 2560.1428
$ws.Range("L3").Value = 2821.3
$ws.Range("M3").Value = -2446.1428
$ws.Range("N3").Value = -3049.3
$ws.Range("H86").Value = 2168.8276
$ws.Range("I86").Value = 1654.8
$ws.Range("J86").Value = 3311.111
$ws.Range("K86").Value = 1654.8
$ws.Range("L86").Value = 3311.111
$ws.Range("M86").Value = -531.8
$ws.Range("N86").Value = -5557.111
$ws.Range("H89").Value = 2168.8276
$ws.Range("I89").Value = 1654.8
$ws.Range("J89").Value = 3311.111
$ws.Range("K89").Value = 8274
$ws.Range("L89").Value = 16555.555
$ws.Range("M89").Value = -2658
$ws.Range("N89").Value = -27787.555
$ws.Range("H134").Value = 52162.855
$ws.Range("I134").Value = 3758.0908
$ws.Range("J134").Value = 108210.48
$ws.Range("K134").Value = 11274.2724
$ws.Range("L134").Value = 324631.44
$ws.Range("M134").Value = -8739.2724
$ws.Range("N134").Value = -329701.44

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4365.967
$ws.Range("I58").Value = 760.5294
$ws.Range("K58").Value = 760.5294
$ws.Range("M58").Value = -557.5294
$ws.Range("H132").Value = 2597.1724
$ws.Range("I132").Value = 1873.4762
$ws.Range("J132").Value = 4496.875
$ws.Range("K132").Value = 5620.4286
$ws.Range("L132").Value = 13490.625
$ws.Range("M132").Value = -3090.4286
$ws.Range("N132").Value = -18550.625
$ws.Range("H134").Value = 3662.4443
$ws.Range("J134").Value = 4199.857
$ws.Range("L134").Value = 12599.571
$ws.Range("N134").Value = -17669.571
$ws.Range("H136").Value = 4365.967
$ws.Range("I136").Value = 760.5294
$ws.Range("K136").Value = 2281.5882
$ws.Range("M136").Value = 268.4117999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1431.8334
$ws.Range("J5").Value = 2935
$ws.Range("L5").Value = 8805
$ws.Range("N5").Value = -9029
$ws.Range("H113").Value = 573.4091
$ws.Range("J113").Value = 583.2432
$ws.Range("L113").Value = 1749.7296
$ws.Range("N113").Value = -6089.729600000001
$ws.Range("H135").Value = 1431.8334
$ws.Range("J135").Value = 2935
$ws.Range("L135").Value = 26415
$ws.Range("N135").Value = -31485

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4521.1904
$ws.Range("I70").Value = 4071.3572
$ws.Range("J70").Value = 5420.857
$ws.Range("K70").Value = 4071.3572
$ws.Range("L70").Value = 5420.857
$ws.Range("M70").Value = -3801.3572
$ws.Range("N70").Value = -5960.857
$ws.Range("H73").Value = 4521.1904
$ws.Range("I73").Value = 4071.3572
$ws.Range("J73").Value = 5420.857
$ws.Range("K73").Value = 4071.3572
$ws.Range("L73").Value = 5420.857
$ws.Range("M73").Value = -3135.3572
$ws.Range("N73").Value = -7292.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4100
$ws.Range("I68").Value = 4250
$ws.Range("J68").Value = 3800
$ws.Range("K68").Value = 4250
$ws.Range("L68").Value = 3800
$ws.Range("M68").Value = -3501
$ws.Range("N68").Value = -5298
$ws.Range("H70").Value = 29900
$ws.Range("J70").Value = 29900
$ws.Range("L70").Value = 29900
$ws.Range("N70").Value = -30440
$ws.Range("H71").Value = 4100
$ws.Range("I71").Value = 4250
$ws.Range("J71").Value = 3800
$ws.Range("K71").Value = 21250
$ws.Range("L71").Value = 19000
$ws.Range("M71").Value = -17506
$ws.Range("N71").Value = -26488
$ws.Range("H73").Value = 29900
$ws.Range("J73").Value = 29900
$ws.Range("L73").Value = 29900
$ws.Range("N73").Value = -31772
